$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns added alongside the existing ORG_GAOL_PK / NAME / ABBR_NAME
# headers: ORG_GAOL_IDENOLD, ORG_GAOL_IDENNEW, ORG_GAOL_STATUS.
$ws.Range("D1").Value = "ORG_GAOL_IDENOLD"
$ws.Range("E1").Value = "ORG_GAOL_IDENNEW"
$ws.Range("F1").Value = "ORG_GAOL_STATUS"

# Match the same left-aligned header style already used by A1:C1.
$ws.Range("D1:F1").HorizontalAlignment = -4131

# Leave the active selection where the author's saved view left it.
$ws.Range("E6").Select()
